# Fix spelling/accent issues ("tildes") across the MATRIZ and DICCIONARIO
# sheets, plus a couple of small wording corrections, matching the
# "Add files via upload" re-upload of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "MATRIZ" -------------------------------------------------------
$matriz = $wb.Worksheets.Item("MATRIZ")

# Row 3 - Corrupcion del codigo actual
$matriz.Cells.Item(3, 2).Value = "Corrupción del código actual"
$matriz.Cells.Item(3, 5).Value = "Crítica"
$matriz.Cells.Item(3, 6).Value = "Siempre se trabaja sobre una copia del proyecto original, y solamente se implementarán cambios si es seguro."

# Row 4 - Falta de personal
$matriz.Cells.Item(4, 5).Value = "Crítica"
$matriz.Cells.Item(4, 6).Value = "Se eliminarán características que no sean totalmente críticas."

# Row 5 - Falta de caracteristicas
$matriz.Cells.Item(5, 2).Value = "Falta de características"
$matriz.Cells.Item(5, 6).Value = "Se implementaran las características faltantes con el tiempo sobrante."

# Row 6 - Documentacion con falta de detalle
$matriz.Cells.Item(6, 2).Value = "Documentación con falta de detalle"

# Row 7 - Falta de requerimientos
$matriz.Cells.Item(7, 6).Value = "Se agregaran los requerimientos faltas solamente si son críticamente necesarios, del otro caso serán omitidos."

# Row 8 - Base de datos erronea
$matriz.Cells.Item(8, 2).Value = "Base de datos errónea "
$matriz.Cells.Item(8, 5).Value = "Catastrófica"
$matriz.Cells.Item(8, 6).Value = "La aplicación se pondrá en cuarentena momentáneamente para solucionar el problema lo antes posible"

# Row 9 - Falta de equipo de desarrollo
$matriz.Cells.Item(9, 3).Value = "Tecnico "
$matriz.Cells.Item(9, 6).Value = "Se modifica el cronograma para sustentar el tiempo perdido a causa del equipo faltante"

# --- Sheet "DICCIONARIO" ---------------------------------------------------
$diccionario = $wb.Worksheets.Item("DICCIONARIO")

# Row 6 - Catastrofica definition
$diccionario.Cells.Item(6, 2).Value = "Catastrófica: Significa el fracaso del proyecto o afectaría gravemente la ejecución/continuación/término del mismo"
